# "added theorems to paper"
#
# Inserts a new "Title and Content" slide right before the existing
# "Related Work" slide (which was slide 4, becomes slide 5) and fills in
# the three SAG-re theorem statements in its body placeholder. The title
# placeholder is left blank, matching the authored slide.

$p = $ppt.ActivePresentation

# "Related Work" is slide 4 today; insert the new slide in front of it so
# it becomes the new slide 4 and "Related Work" slides down to 5.
# Layout 2 == ppLayoutText ("Title and Content"), same layout used by the
# other content slides in this deck.
$newSlide = $p.Slides.Add(4, 2)

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Theorem 1: SAG-re has convergence rate at least as fast as the original SAG.`rTheorem 2: Despite re-computation, SAG-re has asymptotic time complexity as efficient as any gradient method having the lowest iteration cost, namely stochastic gradient.`rTheorem 3: Despite storing the memory gradients, SAG-re has asymptotic space complexity as compact as memory-less gradient methods."
